$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)
$fs = $s.ThemeColorScheme
Write-Output "ok"
# check font scheme access
try {
  $m = $p.SlideMaster
  $ts = $m.TextStyles
  Write-Output "textstyles count: $($ts.Count)"
} catch {
  Write-Output "ERR: $_"
}
